$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (D=Price, E=Volume(1h)).
# D-column values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the original inlineStr cells) instead
# of auto-converting number-like strings (e.g. "102.20" -> 102.2).
$ws.Range("D2").Value = "'27.084.43"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "'1.849.34"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'309.59"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +2.68%  "
$ws.Range("D8").Value = "'0.3685"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "'0.07240"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").Value = "'0.9331"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'19.89"
$ws.Range("D12").Value = "'0.07785"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'1.838.73"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'5.394"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "'6.488"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "'89.16"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'0.000008688"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'1.013"
$ws.Range("D20").Value = "'27.134.51"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("D22").Value = "'5.053"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'1.938"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "'153.03"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "'18.37"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'1.987"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "'114.65"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Value = "'4.900"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "'0.08863"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("D32").Value = "'1.180"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.518"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7394"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'2.691"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").Value = "'1.114"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").Value = "'0.01976"
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("D38").Value = "'0.05267"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "'2.968"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.5268"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'7.034"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "'0.1526"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "'8.289"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'10.55"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'0.4739"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'102.20"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "'66.02"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "'0.06057"
$ws.Range("D51").Value = "'0.8935"
$ws.Range("E51").Value = "  +3.86%  "

# Row 33/34: ImmutableX and Filecoin swapped positions in the ranking.
# (B/C/D/E for rows 33 and 34 already updated above.)
